$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell-level updates per the revised loc combination results
$ws.Range("B1").Value2 = 0.524
$ws.Range("D1").Value2 = 0.971
$ws.Range("E1").Value2 = 0.267
$ws.Range("B2").Value2 = 0.882
$ws.Range("B3").Value2 = 0.86
$ws.Range("B4").Value2 = 0.55
$ws.Range("D4").Value2 = 0.962
$ws.Range("E4").Value2 = 0.4
$ws.Range("B5").Value2 = 0.533
$ws.Range("D5").Value2 = 0.984
$ws.Range("E5").Value2 = 0.267
$ws.Range("B6").Value2 = 1.563
$ws.Range("C6").Value2 = 0.005
$ws.Range("D6").Value2 = 0.972
$ws.Range("E6").Value2 = 0.567
$ws.Range("F6").Value2 = "Muhammad Iqbal Baqi"
$ws.Range("G6").Value2 = "Benar"
$ws.Range("B7").Value2 = 1.167
$ws.Range("C7").Value2 = 0.004
$ws.Range("D7").Value2 = 0.965
$ws.Range("E7").Value2 = 0.467
$ws.Range("B8").Value2 = 1.701
$ws.Range("C8").Value2 = 0.006
$ws.Range("D8").Value2 = 0.995
$ws.Range("E8").Value2 = 0.933
$ws.Range("B9").Value2 = 1.164
$ws.Range("C9").Value2 = 0.004
$ws.Range("D9").Value2 = 0.979
$ws.Range("E9").Value2 = 0.6
$ws.Range("F9").Value2 = "Muhammad Iqbal Baqi"
$ws.Range("G9").Value2 = "Benar"
$ws.Range("B10").Value2 = 1.51
$ws.Range("C10").Value2 = 0.005
$ws.Range("D10").Value2 = 0.977
$ws.Range("E10").Value2 = 0.567
$ws.Range("B11").Value2 = 0.87
$ws.Range("D11").Value2 = 0.971
$ws.Range("B12").Value2 = 1.117
$ws.Range("C12").Value2 = 0.004
$ws.Range("B13").Value2 = 0.912
$ws.Range("D13").Value2 = 0.967
$ws.Range("E13").Value2 = 0.567
$ws.Range("B14").Value2 = 0.952
$ws.Range("D14").Value2 = 0.942
$ws.Range("E14").Value2 = 0.367
$ws.Range("B15").Value2 = 0.956
$ws.Range("D15").Value2 = 0.983
$ws.Range("E15").Value2 = 0.3
$ws.Range("B16").Value2 = 0.762
$ws.Range("C16").Value2 = 0.003
$ws.Range("D16").Value2 = 0.985
$ws.Range("E16").Value2 = 0.333
$ws.Range("B17").Value2 = 0.752
$ws.Range("D17").Value2 = 0.958
$ws.Range("E17").Value2 = 0.467
$ws.Range("B18").Value2 = 1.07
$ws.Range("C18").Value2 = 0.004
$ws.Range("D18").Value2 = 0.985
$ws.Range("E18").Value2 = 0.467
$ws.Range("B19").Value2 = 0.68
$ws.Range("D19").Value2 = 0.97
$ws.Range("B20").Value2 = 1.006
$ws.Range("D20").Value2 = 0.96
$ws.Range("B21").Value2 = 0.585
$ws.Range("D21").Value2 = 0.978
$ws.Range("E21").Value2 = 0.233
$ws.Range("B22").Value2 = 1.312
$ws.Range("D22").Value2 = 0.967
$ws.Range("E22").Value2 = 0.233
$ws.Range("B23").Value2 = 1.187
$ws.Range("C23").Value2 = 0.004
$ws.Range("B24").Value2 = 1.432
$ws.Range("C24").Value2 = 0.005
$ws.Range("D24").Value2 = 0.978
$ws.Range("E24").Value2 = 0.333
$ws.Range("B25").Value2 = 1.372
$ws.Range("C25").Value2 = 0.005
$ws.Range("D25").Value2 = 0.981
$ws.Range("E25").Value2 = 0.333
$ws.Range("B26").Value2 = 1.117
$ws.Range("C26").Value2 = 0.004
$ws.Range("D26").Value2 = 0.895
$ws.Range("E26").Value2 = 0.3
$ws.Range("B27").Value2 = 1.412
$ws.Range("C27").Value2 = 0.005
$ws.Range("D27").Value2 = 0.986
$ws.Range("B28").Value2 = 1.425
$ws.Range("C28").Value2 = 0.005
$ws.Range("D28").Value2 = 0.98
$ws.Range("B29").Value2 = 0.886
$ws.Range("B30").Value2 = 1.556
$ws.Range("C30").Value2 = 0.005
$ws.Range("D30").Value2 = 0.979
$ws.Range("E30").Value2 = 0.7
$ws.Range("B31").Value2 = 0.748
$ws.Range("D31").Value2 = 0.981
$ws.Range("E31").Value2 = 0.267
$ws.Range("B32").Value2 = 0.762
$ws.Range("C32").Value2 = 0.003
$ws.Range("D32").Value2 = 0.981
$ws.Range("B33").Value2 = 1.11
$ws.Range("C33").Value2 = 0.004
$ws.Range("D33").Value2 = 0.964
$ws.Range("E33").Value2 = 0.233
$ws.Range("B34").Value2 = 1.58
$ws.Range("C34").Value2 = 0.005
$ws.Range("D34").Value2 = 0.978
$ws.Range("E34").Value2 = 0.2
$ws.Range("B35").Value2 = 1.448
$ws.Range("C35").Value2 = 0.005
$ws.Range("D35").Value2 = 0.984
$ws.Range("E35").Value2 = 0.267
$ws.Range("F35").Value2 = "Tidak Diketahui"
$ws.Range("G35").Value2 = "Salah"
$ws.Range("A36").Value2 = "FY_4.png"
$ws.Range("B36").Value2 = 1.268
$ws.Range("C36").Value2 = 0.004
$ws.Range("D36").Value2 = 0.964
$ws.Range("E36").Value2 = 0.367
$ws.Range("A37").Value2 = "TO_1.png"
$ws.Range("B37").Value2 = 0.809
$ws.Range("D37").Value2 = 0.962
$ws.Range("E37").Value2 = 0.333
$ws.Range("F37").Value2 = "Tidak Diketahui"
$ws.Range("G37").Value2 = "Salah"
$ws.Range("A38").Value2 = "TO_2.png"
$ws.Range("B38").Value2 = 1.19
$ws.Range("C38").Value2 = 0.004
$ws.Range("D38").Value2 = 0.984
$ws.Range("E38").Value2 = 0.333
$ws.Range("A39").Value2 = "TO_3.png"
$ws.Range("B39").Value2 = 0.828
$ws.Range("C39").Value2 = 0.003
$ws.Range("D39").Value2 = 0.977
$ws.Range("E39").Value2 = 0.4
$ws.Range("A40").Value2 = "TO_4.png"
$ws.Range("B40").Value2 = 3.329
$ws.Range("C40").Value2 = 0.011
$ws.Range("D40").Value2 = 0.836
$ws.Range("E40").Value2 = 0.4
$ws.Range("A41").Value2 = "TO_5.png"
$ws.Range("B41").Value2 = 3.211
$ws.Range("C41").Value2 = 0.01
$ws.Range("D41").Value2 = 0.825
$ws.Range("E41").Value2 = 0.5
$ws.Range("F41").Value2 = "Muhammad Iqbal Baqi"
$ws.Range("G41").Value2 = "Salah"
$ws.Range("A42").Value2 = "TD_1.png"
$ws.Range("B42").Value2 = 2.343
$ws.Range("C42").Value2 = 0.008
$ws.Range("D42").Value2 = 0.849
$ws.Range("E42").Value2 = 0.6
$ws.Range("F42").Value2 = "Muhammad Iqbal Baqi"
$ws.Range("G42").Value2 = "Salah"
$ws.Range("A43").Value2 = "TD_2.png"
$ws.Range("B43").Value2 = 2.529
$ws.Range("C43").Value2 = 0.008
$ws.Range("D43").Value2 = 0.839
$ws.Range("E43").Value2 = 0.333
$ws.Range("A44").Value2 = "TD_3.png"
$ws.Range("B44").Value2 = 1.057
$ws.Range("C44").Value2 = 0.004
$ws.Range("D44").Value2 = 0.941
$ws.Range("A45").Value2 = "TD_4.png"
$ws.Range("B45").Value2 = 1.224
$ws.Range("C45").Value2 = 0.004
$ws.Range("D45").Value2 = 0.979
$ws.Range("E45").Value2 = 0.233
$ws.Range("F45").Value2 = "Tidak Diketahui"
$ws.Range("G45").Value2 = "Benar"
